# Updates the correlation-factor table (Table 9) with revised "ABUNDANCE"
# values for COD (column 2) and HAKE (column 6) per country row, reflecting
# the new correlation results referenced in the commit message.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellRange = $t.Cell(3, 2).Range
$cellRange.Find.Execute("0.092", $true, $false, $false, $false, $false, $true, 1, $false, "0.138", 2) | Out-Null
$cellRange = $t.Cell(3, 6).Range
$cellRange.Find.Execute("0.764", $true, $false, $false, $false, $false, $true, 1, $false, "0.709", 2) | Out-Null
$cellRange = $t.Cell(4, 2).Range
$cellRange.Find.Execute("0.387", $true, $false, $false, $false, $false, $true, 1, $false, "0.398", 2) | Out-Null
$cellRange = $t.Cell(4, 6).Range
$cellRange.Find.Execute("0.764", $true, $false, $false, $false, $false, $true, 1, $false, "0.709", 2) | Out-Null
$cellRange = $t.Cell(5, 2).Range
$cellRange.Find.Execute("0.383", $true, $false, $false, $false, $false, $true, 1, $false, "0.363", 2) | Out-Null
$cellRange = $t.Cell(5, 6).Range
$cellRange.Find.Execute("0.764", $true, $false, $false, $false, $false, $true, 1, $false, "0.709", 2) | Out-Null
$cellRange = $t.Cell(6, 2).Range
$cellRange.Find.Execute("0.379", $true, $false, $false, $false, $false, $true, 1, $false, "0.357", 2) | Out-Null
$cellRange = $t.Cell(7, 2).Range
$cellRange.Find.Execute("0.309", $true, $false, $false, $false, $false, $true, 1, $false, "0.339", 2) | Out-Null
$cellRange = $t.Cell(8, 2).Range
$cellRange.Find.Execute("0.470", $true, $false, $false, $false, $false, $true, 1, $false, "0.466", 2) | Out-Null
$cellRange = $t.Cell(8, 6).Range
$cellRange.Find.Execute("0.595", $true, $false, $false, $false, $false, $true, 1, $false, "0.425", 2) | Out-Null
$cellRange = $t.Cell(9, 2).Range
$cellRange.Find.Execute("0.329", $true, $false, $false, $false, $false, $true, 1, $false, "0.368", 2) | Out-Null
$cellRange = $t.Cell(9, 6).Range
$cellRange.Find.Execute("0.595", $true, $false, $false, $false, $false, $true, 1, $false, "0.425", 2) | Out-Null
$cellRange = $t.Cell(10, 2).Range
$cellRange.Find.Execute("0.199", $true, $false, $false, $false, $false, $true, 1, $false, "0.139", 2) | Out-Null
$cellRange = $t.Cell(11, 2).Range
$cellRange.Find.Execute("0.199", $true, $false, $false, $false, $false, $true, 1, $false, "0.139", 2) | Out-Null
$cellRange = $t.Cell(12, 2).Range
$cellRange.Find.Execute("0.408", $true, $false, $false, $false, $false, $true, 1, $false, "0.520", 2) | Out-Null
$cellRange = $t.Cell(12, 6).Range
$cellRange.Find.Execute("0.764", $true, $false, $false, $false, $false, $true, 1, $false, "0.709", 2) | Out-Null
$cellRange = $t.Cell(13, 2).Range
$cellRange.Find.Execute("0.379", $true, $false, $false, $false, $false, $true, 1, $false, "0.357", 2) | Out-Null
$cellRange = $t.Cell(14, 2).Range
$cellRange.Find.Execute("0.449", $true, $false, $false, $false, $false, $true, 1, $false, "0.484", 2) | Out-Null
$cellRange = $t.Cell(14, 6).Range
$cellRange.Find.Execute("0.595", $true, $false, $false, $false, $false, $true, 1, $false, "0.425", 2) | Out-Null
$cellRange = $t.Cell(15, 2).Range
$cellRange.Find.Execute("0.199", $true, $false, $false, $false, $false, $true, 1, $false, "0.139", 2) | Out-Null
$cellRange = $t.Cell(16, 2).Range
$cellRange.Find.Execute("0.304", $true, $false, $false, $false, $false, $true, 1, $false, "0.329", 2) | Out-Null
$cellRange = $t.Cell(16, 6).Range
$cellRange.Find.Execute("0.764", $true, $false, $false, $false, $false, $true, 1, $false, "0.709", 2) | Out-Null
